# chore: update Sheets via scheduled runner
# Refresh pricing-derived figures (currentAveragePrice* / LevePrice* / LeveProfit*)
# across the per-job leve profit tables.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 34800
$ws.Range("J68").Value = 34800
$ws.Range("L68").Value = 34800
$ws.Range("N68").Value = -36298

$ws.Range("H71").Value = 34800
$ws.Range("J71").Value = 34800
$ws.Range("L71").Value = 104400
$ws.Range("N71").Value = -111888

$ws.Range("H80").Value = 347.8125
$ws.Range("I80").Value = 282.46155
$ws.Range("J80").Value = 631
$ws.Range("K80").Value = 847.38465
$ws.Range("L80").Value = 1893
$ws.Range("M80").Value = 150.61535
$ws.Range("N80").Value = -3889

$ws.Range("H83").Value = 347.8125
$ws.Range("I83").Value = 282.46155
$ws.Range("J83").Value = 631
$ws.Range("K83").Value = 2542.15395
$ws.Range("L83").Value = 5679
$ws.Range("M83").Value = 2449.84605
$ws.Range("N83").Value = -15663

$ws.Range("H86").Value = 6151.593
$ws.Range("I86").Value = 6653.9375
$ws.Range("J86").Value = 5420.909
$ws.Range("K86").Value = 6653.9375
$ws.Range("L86").Value = 5420.909
$ws.Range("M86").Value = -5530.9375
$ws.Range("N86").Value = -7666.909

$ws.Range("H89").Value = 6151.593
$ws.Range("I89").Value = 6653.9375
$ws.Range("J89").Value = 5420.909
$ws.Range("K89").Value = 33269.6875
$ws.Range("L89").Value = 27104.545
$ws.Range("M89").Value = -27653.6875
$ws.Range("N89").Value = -38336.545

$ws.Range("H94").Value = 2757.7144
$ws.Range("I94").Value = 2757.7144
$ws.Range("K94").Value = 2757.7144
$ws.Range("M94").Value = -2306.7144

$ws.Range("H106").Value = 43479496
$ws.Range("I106").Value = 55556324
$ws.Range("J106").Value = 2920
$ws.Range("K106").Value = 55556324
$ws.Range("L106").Value = 2920
$ws.Range("M106").Value = -55555693
$ws.Range("N106").Value = -4182

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 181.7
$ws.Range("I4").Value = 166.71428
$ws.Range("J4").Value = 216.66667
$ws.Range("K4").Value = 166.71428
$ws.Range("L4").Value = 216.66667
$ws.Range("M4").Value = -50.71428
$ws.Range("N4").Value = -448.66667

$ws.Range("H32").Value = 11241387
$ws.Range("I32").Value = 4186.2
$ws.Range("K32").Value = 4186.2
$ws.Range("M32").Value = -3899.2

$ws.Range("H132").Value = 1341.8959
$ws.Range("I132").Value = 947.7105
$ws.Range("J132").Value = 2839.8
$ws.Range("K132").Value = 2843.1315
$ws.Range("L132").Value = 8519.400000000001
$ws.Range("M132").Value = -313.1315
$ws.Range("N132").Value = -13579.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 481.1875
$ws.Range("I94").Value = 370.8095
$ws.Range("K94").Value = 370.8095
$ws.Range("M94").Value = 80.19049999999999

$ws.Range("H132").Value = 49490
$ws.Range("J132").Value = 49490
$ws.Range("L132").Value = 49490
$ws.Range("N132").Value = -59610

$ws.Range("H134").Value = 3179130.8
$ws.Range("I134").Value = 980.913
$ws.Range("J134").Value = 9270584
$ws.Range("K134").Value = 2942.739
$ws.Range("L134").Value = 27811752
$ws.Range("M134").Value = -407.739
$ws.Range("N134").Value = -27816822

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 50102420
$ws.Range("J86").Value = 74314.28999999999
$ws.Range("L86").Value = 74314.28999999999
$ws.Range("N86").Value = -76560.28999999999

$ws.Range("H89").Value = 50102420
$ws.Range("J89").Value = 74314.28999999999
$ws.Range("L89").Value = 371571.45
$ws.Range("N89").Value = -382803.45

$ws.Range("H105").Value = 6086.8
$ws.Range("I105").Value = 7150.4
$ws.Range("J105").Value = 2896
$ws.Range("K105").Value = 7150.4
$ws.Range("L105").Value = 2896
$ws.Range("M105").Value = -5403.4
$ws.Range("N105").Value = -6390

$ws.Range("H132").Value = 18520254
$ws.Range("I132").Value = 1270.5
$ws.Range("J132").Value = 55558224
$ws.Range("K132").Value = 3811.5
$ws.Range("L132").Value = 166674672
$ws.Range("M132").Value = -1281.5
$ws.Range("N132").Value = -166679732

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 112.4
$ws.Range("I6").Value = 67
$ws.Range("J6").Value = 294
$ws.Range("K6").Value = 201
$ws.Range("L6").Value = 882
$ws.Range("M6").Value = -88
$ws.Range("N6").Value = -1108

$ws.Range("H11").Value = 18843.75
$ws.Range("I11").Value = 124
$ws.Range("J11").Value = 50043.332
$ws.Range("K11").Value = 372
$ws.Range("L11").Value = 150129.996
$ws.Range("M11").Value = -232
$ws.Range("N11").Value = -150409.996

$ws.Range("H137").Value = 500000000
$ws.Range("I137").Value = 500000000
$ws.Range("K137").Value = 1500000000
$ws.Range("M137").Value = -1499994900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6033.636
$ws.Range("I132").Value = 1196.9474
$ws.Range("J132").Value = 36666
$ws.Range("K132").Value = 3590.8422
$ws.Range("L132").Value = 109998
$ws.Range("M132").Value = -1060.8422
$ws.Range("N132").Value = -115058

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1148.2632
$ws.Range("I93").Value = 1059
$ws.Range("K93").Value = 1059
$ws.Range("M93").Value = 189

$ws.Range("H132").Value = 12971.479
$ws.Range("I132").Value = 3948.4443
$ws.Range("J132").Value = 18772
$ws.Range("K132").Value = 11845.3329
$ws.Range("L132").Value = 56316
$ws.Range("M132").Value = -9315.332900000001
$ws.Range("N132").Value = -61376

$ws.Range("H136").Value = 4468
$ws.Range("I136").Value = 4038.9583
$ws.Range("J136").Value = 6184.1665
$ws.Range("K136").Value = 12116.8749
$ws.Range("L136").Value = 18552.4995
$ws.Range("M136").Value = -9566.874899999999
$ws.Range("N136").Value = -23652.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 8354
$ws.Range("I38").Value = 10000
$ws.Range("J38").Value = 7531
$ws.Range("K38").Value = 10000
$ws.Range("L38").Value = 7531
$ws.Range("M38").Value = -9527
$ws.Range("N38").Value = -8477

$ws.Range("H75").Value = 19899.75
$ws.Range("J75").Value = 19899.75
$ws.Range("L75").Value = 19899.75
$ws.Range("N75").Value = -21771.75

$ws.Range("H78").Value = 19899.75
$ws.Range("J78").Value = 19899.75
$ws.Range("L78").Value = 59699.25
$ws.Range("N78").Value = -69059.25

$ws.Range("H132").Value = 40388.656
$ws.Range("I132").Value = 52660.15
$ws.Range("J132").Value = 13118.667
$ws.Range("K132").Value = 157980.45
$ws.Range("L132").Value = 39356.001
$ws.Range("M132").Value = -155450.45
$ws.Range("N132").Value = -44416.001
